$d = $word.ActiveDocument

# Remove the three paragraphs that follow the
# "LOQ4087: ... (Requisito fraco)" paragraph:
#   - the blank paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: ... Creative Commons Attribution"
# The blank paragraph that originally followed those (right before the
# page-break paragraph) is left in place.

$start = $null
$end = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*LOQ4087*") {
        $start = $p.Next().Range.Start
    }
    if ($t -like "*© 2020*") {
        $end = $p.Range.End
    }
}

if ($start -ne $null -and $end -ne $null) {
    $d.Range($start, $end).Delete()
}
